# Auto-generated edit script applying cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Mon Sep  9 08:41:27 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.858.40"
$ws.Range("E2").Value = "'  +0.96%  "

$ws.Range("D3").Value = "'2.293.04"
$ws.Range("E3").Value = "'  +0.48%  "

$ws.Range("E4").Value = "'  +0.13%  "

$ws.Range("D5").Value = "'507.29"
$ws.Range("E5").Value = "'  +0.41%  "

$ws.Range("D6").Value = "'129.72"
$ws.Range("E6").Value = "'  +0.14%  "

$ws.Range("E7").Value = "'  -0.30%  "

$ws.Range("E8").Value = "'  +0.43%  "

$ws.Range("D9").Value = "'2.316.96"
$ws.Range("E9").Value = "'  +1.18%  "

$ws.Range("D10").Value = "'0.0972"
$ws.Range("E10").Value = "'  +1.61%  "

$ws.Range("E11").Value = "'  +1.72%  "

$ws.Range("B12").Value = "'Toncoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'5.06"
$ws.Range("E12").Value = "'  +7.05%  "

$ws.Range("B13").Value = "'Cardano"
$ws.Range("C13").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.340"
$ws.Range("E13").Value = "'  +1.93%  "

$ws.Range("D14").Value = "'23.89"
$ws.Range("E14").Value = "'  +4.10%  "

$ws.Range("D15").Value = "'2.706.09"
$ws.Range("E15").Value = "'  +0.68%  "

$ws.Range("D16").Value = "'54.883.28"

$ws.Range("E17").Value = "'  +1.58%  "

$ws.Range("D18").Value = "'2.272.48"
$ws.Range("E18").Value = "'  -0.56%  "

$ws.Range("D19").Value = "'10.69"
$ws.Range("E19").Value = "'  +3.83%  "

$ws.Range("E20").Value = "'  +1.49%  "

$ws.Range("D22").Value = "'309.78"
$ws.Range("E22").Value = "'  +1.48%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "'  -0.28%  "

$ws.Range("D24").Value = "'60.59"
$ws.Range("E24").Value = "'  -2.00%  "

$ws.Range("D25").Value = "'0.991"
$ws.Range("E25").Value = "'  -0.83%  "

$ws.Range("E26").Value = "'  +0.20%  "

$ws.Range("E27").Value = "'  +2.49%  "

$ws.Range("D28").Value = "'172.67"
$ws.Range("E28").Value = "'  -1.09%  "

$ws.Range("D29").Value = "'6.15"
$ws.Range("E29").Value = "'  +2.06%  "

$ws.Range("D30").Value = "'0.0₃0708"
$ws.Range("E30").Value = "'  +2.46%  "

$ws.Range("E31").Value = "'  +0.56%  "

$ws.Range("E32").Value = "'  +5.39%  "

$ws.Range("B33").Value = "'EthereumClassic"
$ws.Range("C33").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.11"
$ws.Range("E33").Value = "'  +1.76%  "

$ws.Range("B34").Value = "'USDe"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "'  -0.02%  "

$ws.Range("E35").Value = "'  -0.27%  "

$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.23"
$ws.Range("E36").Value = "'  +2.63%  "

$ws.Range("B37").Value = "'SuiNetwork"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "'0.916"
$ws.Range("E37").Value = "'  -5.21%  "

$ws.Range("D38").Value = "'3.88"
$ws.Range("E38").Value = "'  +3.04%  "

$ws.Range("D39").Value = "'36.76"
$ws.Range("E39").Value = "'  +1.96%  "

$ws.Range("E40").Value = "'  +1.24%  "

$ws.Range("E41").Value = "'  +1.91%  "

$ws.Range("D42").Value = "'135.43"
$ws.Range("E42").Value = "'  +8.17%  "

$ws.Range("B43").Value = "'RenderToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.11"
$ws.Range("E43").Value = "'  +5.30%  "

$ws.Range("B44").Value = "'Filecoin"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.44"
$ws.Range("E44").Value = "'  +1.62%  "

$ws.Range("D45").Value = "'257.34"
$ws.Range("E45").Value = "'  +6.04%  "

$ws.Range("D46").Value = "'0.0504"
$ws.Range("E46").Value = "'  +1.56%  "

$ws.Range("E47").Value = "'  +2.08%  "

$ws.Range("E48").Value = "'  +1.11%  "

$ws.Range("E49").Value = "'  +1.42%  "

$ws.Range("D50").Value = "'0.0209"
$ws.Range("E50").Value = "'  +1.30%  "

$ws.Range("E51").Value = "'  +0.36%  "

